# Apply the edit described by the diff:
#  - Remove the "крупа" data row (originally row 2), which also removes its
#    shared string and shifts all subsequent rows up by one.
#  - Remove the "чича" and "fertgreg" data rows (originally the last two rows),
#    removing their shared strings as well.
#
# Resulting data (after deletions and the natural row shift caused by them):
#   Name   | Priority
#   лук    | 5
#   чеснок | 6
#   гречка | 4
#   лапша  | 1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two trailing rows first (originally rows 7 and 8, i.e. "чича" and
# "fertgreg"), so row indices for the earlier rows remain stable while we do
# this. Then delete row 2 ("крупа") which shifts everything below it up.
$ws.Rows.Item(8).EntireRow.Delete() | Out-Null
$ws.Rows.Item(7).EntireRow.Delete() | Out-Null
$ws.Rows.Item(2).EntireRow.Delete() | Out-Null
